$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the Excel Table (ListObject) that holds the reponedor data
$tbl = $ws.ListObjects.Item("Tabla1")

# Add two new rows at the end of the table for the new local "Jumbo El Llano"
$row1 = $tbl.ListRows.Add()
$row2 = $tbl.ListRows.Add()

$row1.Range.Item(1, 1).Value = "Jumbo El Llano"
$row1.Range.Item(1, 2).Value = "J513"
$row1.Range.Item(1, 3).Value = "J513"
$row1.Range.Item(1, 4).Value = 4
$row1.Range.Item(1, 5).Value = "Antonio"
$row1.Range.Item(1, 6).Value = 5
$row1.Range.Item(1, 7).Value = "viernes "

$row2.Range.Item(1, 1).Value = "Jumbo El Llano"
$row2.Range.Item(1, 2).Value = "J513"
$row2.Range.Item(1, 3).Value = "J513"
$row2.Range.Item(1, 4).Value = 4
$row2.Range.Item(1, 5).Value = "Antonio"
$row2.Range.Item(1, 6).Value = 6
$row2.Range.Item(1, 7).Value = "sábado "

$ws.Range("A152").Select()
